# "muda de pygame para flask":
#   - header row (1) loses the bold/boxed/centered formatting -> back to default "Normal" style
#   - a handful of point values are corrected on rows 5 and 28
#   - five new survey rows (31-35) are appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the bold/border/center-top header formatting back to the default style ---
$ws.Range("A1:BA1").Style = "Normal"

# --- 2. Point fixes in row 5 ---
$row5Fixes = @{
  "D" = 14.75
  "F" = 9.050000000000001
  "G" = -12.85
  "I" = 12.85
  "M" = 12.85
  "N" = 14.75
}
foreach ($col in $row5Fixes.Keys) {
  $ws.Range($col + "5").Value2 = $row5Fixes[$col]
}

# --- 3. Point fixes in row 28 ---
$row28Fixes = @{
  "B" = 1.9
  "C" = -0.95
  "D" = -2.85
  "G" = -0.95
  "J" = -2.85
  "N" = 2.85
  "R" = 1.5
  "U" = 0.95
  "AF" = -0.95
  "AG" = 0.95
  "AI" = 0.75
}
foreach ($col in $row28Fixes.Keys) {
  $ws.Range($col + "28").Value2 = $row28Fixes[$col]
}

# --- 4. Append new rows 31-35 ---
$newRows = @{
  31 = @("Bruna (T26)", 3.8, 4.75, -4.75, -3.8, 0, 0, 0, 0, -4.75, -3.8, -3.8, -3.8, -3.8, -3.8, -3.8, -3.8, 3.75, -3.8, -3.8, 4.75, -3.8, -3.8, -3.8, -3.8, -3.8, -4.75, 3.8, -3.8, -3.8, -3.8, -3.8, 3.8, -3.8, 0, 0, -3.8, -3.8, -3.8, 3.8, -3.8, -3.8, 3.8, -3.8, -3.8, -3.8, -3.8, 0, 0, 0, 0, -3.8, -3.8)
  32 = @("Ricardo (T26)", 1.9, -1.9, -2.85, -1.9, -1.9, -1.9, -1.9, 0, -1.9, -1.9, -1.9, -1.9, -1.9, -1.9, -1.9, -1.9, 2.25, -1.9, -1.9, 2.85, -1.9, -1.9, -1.9, -1.9, -1.9, -1.9, 1.9, 0, -1.9, -2.85, -1.9, -1.9, 1.9, -0.95, 0, -1.9, -1.15, -1.9, -1.9, 2.85, -1.9, 1.9, 1.9, -1.9, 1.5, -1.9, 0, 0, 0, 0, -1.9, -1.9)
  33 = @("teste", 0, 0, 0.95, 0, 0.95, 0, 0, 0, 0, 0, 0, 0, 0.95, 0.95, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  34 = @("a", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
  35 = @("aaa", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($rowNum in $newRows.Keys) {
  $rowValues = $newRows[$rowNum]
  for ($colIdx = 0; $colIdx -lt $rowValues.Length; $colIdx++) {
    $ws.Cells.Item([int]$rowNum, $colIdx + 1).Value2 = $rowValues[$colIdx]
  }
}
